# Apply edit: fix "surplus number" values in column K (from 1.0565 to 1)
# for the rows that reference the surplus/markup multiplier, and update
# the active selection to E3 (supports longer quotes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer Quote")

# Correct the erroneous surplus number (1.0565 -> 1) in column K
$rows = @(16, 17, 19, 22, 25, 26, 30)
foreach ($r in $rows) {
    $ws.Range("K$r").Value = 1
}

# Move/update the active selection on the sheet to support longer quotes
$ws.Activate()
$ws.Range("E3").Select()
